$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 319.0909
$ws.Range("I19").Value = 216.25
$ws.Range("J19").Value = 377.85715
$ws.Range("K19").Value = 216.25
$ws.Range("L19").Value = 377.85715
$ws.Range("M19").Value = -41.25
$ws.Range("N19").Value = -727.85715
$ws.Range("H62").Value = 23815828
$ws.Range("I62").Value = 1592.9231
$ws.Range("J62").Value = 62513960
$ws.Range("K62").Value = 1592.9231
$ws.Range("L62").Value = 62513960
$ws.Range("M62").Value = -968.9231
$ws.Range("N62").Value = -62515208
$ws.Range("H65").Value = 23815828
$ws.Range("I65").Value = 1592.9231
$ws.Range("J65").Value = 62513960
$ws.Range("K65").Value = 7964.6155
$ws.Range("L65").Value = 312569800
$ws.Range("M65").Value = -4844.6155
$ws.Range("N65").Value = -312576040
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("H100").Value = 23811232
$ws.Range("I100").Value = 1756.7646
$ws.Range("J100").Value = 125001496
$ws.Range("K100").Value = 1756.7646
$ws.Range("L100").Value = 125001496
$ws.Range("M100").Value = -1215.7646
$ws.Range("N100").Value = -125002578
$ws.Range("H123").Value = 50993.332
$ws.Range("J123").Value = 50993.332
$ws.Range("L123").Value = 50993.332
$ws.Range("N123").Value = -60793.332
$ws.Range("H128").Value = 49097.5
$ws.Range("J128").Value = 49097.5
$ws.Range("L128").Value = 49097.5
$ws.Range("N128").Value = -59057.5
$ws.Range("H133").Value = 40833
$ws.Range("J133").Value = 40833
$ws.Range("L133").Value = 40833
$ws.Range("N133").Value = -50953
$ws.Range("H134").Value = 49290
$ws.Range("J134").Value = 49290
$ws.Range("L134").Value = 49290
$ws.Range("N134").Value = -59430
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23975.283
$ws.Range("I32").Value = 6936.1914
$ws.Range("K32").Value = 6936.1914
$ws.Range("M32").Value = -6649.1914
$ws.Range("H63").Value = 4000
$ws.Range("I63").Value = 2800
$ws.Range("K63").Value = 2800
$ws.Range("M63").Value = -2114
$ws.Range("H66").Value = 4000
$ws.Range("I66").Value = 2800
$ws.Range("K66").Value = 14000
$ws.Range("M66").Value = -10568
$ws.Range("H103").Value = 49888
$ws.Range("J103").Value = 49888
$ws.Range("L103").Value = 49888
$ws.Range("N103").Value = -52232
$ws.Range("H109").Value = 31000
$ws.Range("J109").Value = 31000
$ws.Range("L109").Value = 31000
$ws.Range("N109").Value = -33774
$ws.Range("H123").Value = 48870
$ws.Range("J123").Value = 48870
$ws.Range("L123").Value = 48870
$ws.Range("N123").Value = -58670
$ws.Range("H129").Value = 38551.8
$ws.Range("J129").Value = 38551.8
$ws.Range("L129").Value = 38551.8
$ws.Range("N129").Value = -48551.8
$ws.Range("H130").Value = 48298.168
$ws.Range("J130").Value = 48298.168
$ws.Range("L130").Value = 48298.168
$ws.Range("N130").Value = -58338.168
$ws.Range("H132").Value = 2510.2683
$ws.Range("I132").Value = 2155.4285
$ws.Range("J132").Value = 3274.5386
$ws.Range("K132").Value = 6466.2855
$ws.Range("L132").Value = 9823.6158
$ws.Range("M132").Value = -3936.2855
$ws.Range("N132").Value = -14883.6158
$ws.Range("H133").Value = 29248.715
$ws.Range("J133").Value = 29248.715
$ws.Range("L133").Value = 29248.715
$ws.Range("N133").Value = -34308.715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 45552.5
$ws.Range("J122").Value = 45552.5
$ws.Range("L122").Value = 45552.5
$ws.Range("N122").Value = -55352.5
$ws.Range("H126").Value = 31367.5
$ws.Range("J126").Value = 31367.5
$ws.Range("L126").Value = 31367.5
$ws.Range("N126").Value = -41247.5
$ws.Range("H130").Value = 69980
$ws.Range("J130").Value = 69980
$ws.Range("L130").Value = 69980
$ws.Range("N130").Value = -80020
$ws.Range("H132").Value = 37212
$ws.Range("J132").Value = 37212
$ws.Range("L132").Value = 37212
$ws.Range("N132").Value = -47332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 46074.145
$ws.Range("J20").Value = 46074.145
$ws.Range("L20").Value = 46074.145
$ws.Range("N20").Value = -46546.145
$ws.Range("H30").Value = 46074.145
$ws.Range("J30").Value = 46074.145
$ws.Range("L30").Value = 46074.145
$ws.Range("N30").Value = -46256.145
$ws.Range("H58").Value = 1571.579
$ws.Range("I58").Value = 1390.3704
$ws.Range("J58").Value = 2016.3636
$ws.Range("K58").Value = 1390.3704
$ws.Range("L58").Value = 2016.3636
$ws.Range("M58").Value = -1187.3704
$ws.Range("N58").Value = -2422.3636
$ws.Range("H97").Value = 19994.5
$ws.Range("J97").Value = 19994.5
$ws.Range("L97").Value = 19994.5
$ws.Range("N97").Value = -21976.5
$ws.Range("H127").Value = 45340
$ws.Range("J127").Value = 45340
$ws.Range("L127").Value = 45340
$ws.Range("N127").Value = -55260
$ws.Range("H128").Value = 46074.145
$ws.Range("J128").Value = 46074.145
$ws.Range("L128").Value = 46074.145
$ws.Range("N128").Value = -56034.145
$ws.Range("H135").Value = 41657.145
$ws.Range("J135").Value = 41657.145
$ws.Range("L135").Value = 41657.145
$ws.Range("N135").Value = -51797.145
$ws.Range("H136").Value = 1571.579
$ws.Range("I136").Value = 1390.3704
$ws.Range("J136").Value = 2016.3636
$ws.Range("K136").Value = 4171.1112
$ws.Range("L136").Value = 6049.0908
$ws.Range("M136").Value = -1621.1112
$ws.Range("N136").Value = -11149.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7902.5454
$ws.Range("I3").Value = 6434.2856
$ws.Range("J3").Value = 10472
$ws.Range("K3").Value = 19302.8568
$ws.Range("L3").Value = 31416
$ws.Range("M3").Value = -19190.8568
$ws.Range("N3").Value = -31640
$ws.Range("H97").Value = 634.6667
$ws.Range("I97").Value = 850
$ws.Range("J97").Value = 204
$ws.Range("K97").Value = 2550
$ws.Range("L97").Value = 612
$ws.Range("M97").Value = -2054
$ws.Range("N97").Value = -1604

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 14067.223
$ws.Range("J93").Value = 14067.223
$ws.Range("L93").Value = 14067.223
$ws.Range("N93").Value = -17811.223
$ws.Range("H109").Value = 19440.416
$ws.Range("J109").Value = 19440.416
$ws.Range("L109").Value = 19440.416
$ws.Range("N109").Value = -21520.416
$ws.Range("H126").Value = 22629.23
$ws.Range("I126").Value = 5993.3335
$ws.Range("J126").Value = 27620
$ws.Range("K126").Value = 17980.0005
$ws.Range("L126").Value = 82860
$ws.Range("M126").Value = -15510.0005
$ws.Range("N126").Value = -87800
$ws.Range("H133").Value = 28465.834
$ws.Range("J133").Value = 28465.834
$ws.Range("L133").Value = 28465.834
$ws.Range("N133").Value = -38585.834
$ws.Range("H135").Value = 59223
$ws.Range("J135").Value = 59223
$ws.Range("L135").Value = 59223
$ws.Range("N135").Value = -69363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 18374.75
$ws.Range("J59").Value = 18374.75
$ws.Range("L59").Value = 18374.75
$ws.Range("N59").Value = -19682.75
$ws.Range("H63").Value = 29267.5
$ws.Range("J63").Value = 29267.5
$ws.Range("L63").Value = 29267.5
$ws.Range("N63").Value = -30765.5
$ws.Range("H66").Value = 29267.5
$ws.Range("J66").Value = 29267.5
$ws.Range("L66").Value = 87802.5
$ws.Range("N66").Value = -95290.5
$ws.Range("H96").Value = 9890
$ws.Range("J96").Value = 9890
$ws.Range("L96").Value = 9890
$ws.Range("N96").Value = -15382
$ws.Range("H100").Value = 2089.2144
$ws.Range("I100").Value = 1794.4445
$ws.Range("J100").Value = 2619.8
$ws.Range("K100").Value = 1794.4445
$ws.Range("L100").Value = 2619.8
$ws.Range("M100").Value = -1253.4445
$ws.Range("N100").Value = -3701.8
$ws.Range("H108").Value = 24241.166
$ws.Range("J108").Value = 24241.166
$ws.Range("L108").Value = 24241.166
$ws.Range("N108").Value = -31921.166
$ws.Range("H122").Value = 4929.647
$ws.Range("I122").Value = 4444.8887
$ws.Range("J122").Value = 5475
$ws.Range("K122").Value = 13334.6661
$ws.Range("L122").Value = 16425
$ws.Range("M122").Value = -10884.6661
$ws.Range("N122").Value = -21325
$ws.Range("H123").Value = 38911.285
$ws.Range("J123").Value = 38911.285
$ws.Range("L123").Value = 38911.285
$ws.Range("N123").Value = -48711.285
$ws.Range("H127").Value = 42792
$ws.Range("J127").Value = 42792
$ws.Range("L127").Value = 42792
$ws.Range("N127").Value = -52712
$ws.Range("H128").Value = 44398.168
$ws.Range("J128").Value = 44398.168
$ws.Range("L128").Value = 44398.168
$ws.Range("N128").Value = -54358.168
$ws.Range("H130").Value = 42835.57
$ws.Range("J130").Value = 42835.57
$ws.Range("L130").Value = 42835.57
$ws.Range("N130").Value = -52875.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 23994
$ws.Range("J93").Value = 23994
$ws.Range("L93").Value = 23994
$ws.Range("N93").Value = -28986
$ws.Range("H125").Value = 30837.334
$ws.Range("J125").Value = 30837.334
$ws.Range("L125").Value = 30837.334
$ws.Range("N125").Value = -40677.334
$ws.Range("H132").Value = 2017.8889
$ws.Range("I132").Value = 2218.4517
$ws.Range("J132").Value = 1573.7858
$ws.Range("K132").Value = 6655.355100000001
$ws.Range("L132").Value = 4721.357400000001
$ws.Range("M132").Value = -4125.355100000001
$ws.Range("N132").Value = -9781.357400000001
